# Sector coupled implemented correctly
$wb = $excel.ActiveWorkbook

$wsGen  = $wb.Worksheets.Item("Generator")
$wsLink = $wb.Worksheets.Item("Link")

# --- Generator sheet: "diesel" row's bus changes from "bus 1" to "bus 0" ---
$wsGen.Range("C3").Value = "bus 0"

# --- Link sheet: update p_nom_extendable / efficiency / p_nom / efficiency2 values ---
# p_nom_extendable (column E) must become the literal text "False" (not the Excel
# boolean FALSE), matching how the other True/False flags are stored as shared
# strings in this workbook. A direct Value/Formula assignment of "False" gets
# auto-coerced to a real boolean, so instead we build it with TEXT() (forcing a
# string result) and then convert the formula to a static value via copy / paste
# special, which keeps it typed as text.
function Set-TextFalse($cell) {
    $cell.Formula = '=TEXT(0,"\F\a\l\s\e")'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

Set-TextFalse $wsLink.Range("E2")
$wsLink.Range("I2").Value = 100

Set-TextFalse $wsLink.Range("E3")
$wsLink.Range("H3").Value = 0.5
$wsLink.Range("I3").Value = 100
$wsLink.Range("K3").Value = 0.4

$excel.CutCopyMode = 0

# --- Update selections / active sheet to match the new workbook view state ---
# Generator is no longer the active tab; its selection moves to C4.
$wsGen.Range("C4").Select()

# Link becomes the active tab, with its selection on K4.
$wsLink.Activate()
$wsLink.Range("K4").Select()
